# Insert a new bullet point right after the "For each of the listed files..."
# list item, offering the "git add ." shortcut.

$d = $word.ActiveDocument

# Locate the paragraph to insert after via Find (robust against any
# paragraph re-indexing).
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "For each of the listed files, type git add [filename] (this can also be foldername/) and press enter",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorPara = $searchRange.Paragraphs(1)

# Create a brand new paragraph right after it; Word seeds its pPr/rPr
# (style, numbering, run-mark formatting) from the paragraph it follows.
$anchorPara.Range.InsertParagraphAfter()

$newPara = $anchorPara.Next()
$newRange = $newPara.Range

$boldLeadIn = "Or better yet: "
$rest = "when you" + [char]0x2019 + "re sure you want to add all files use " + [char]0x201C + "git add ." + [char]0x201D

$startPos = $newRange.Start

# Type both runs worth of text first (as one plain run) ...
$newRange.Text = $boldLeadIn + $rest

# ... then apply bold + underline only to the "Or better yet: " lead-in so
# it becomes its own distinctly-formatted run, leaving the remainder (and
# the paragraph mark) in the inherited, un-bolded formatting.
$leadInRange = $d.Range($startPos, $startPos + $boldLeadIn.Length)
$leadInRange.Font.Bold = $true
$leadInRange.Font.Underline = 1
